# Project 3 (Exercise 1) clustering results: fill in the Run #1-#3 WCSS
# values for each data file (columns B:D, rows 3-24) now that the runs have
# finished. The Median (E) / Minimum (F) formulas already in the sheet
# recalculate automatically once the inputs are present, clearing the
# #NUM! errors that came from MEDIAN()/MIN() over empty ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$results = @{
    3  = @(22.384620000000002, 22.384620000000002, 22.429427)
    4  = @(9.7972490000000008, 9.7972490000000008, 9.7972490000000008)
    5  = @(1.1214090000000001, 1.1214090000000001, 1.1214090000000001)
    6  = @(1.57145, 0.52234700000000001, 1.415432)
    7  = @(0.55971700000000002, 0.55971700000000002, 0.55971700000000002)
    8  = @(0.51026800000000005, 0.51026800000000005, 0.51026800000000005)
    9  = @(0.47006199999999998, 0.47006199999999998, 5.2599099999999996)
    10 = @(5.2264080000000002, 0.51057200000000003, 2.2968929999999999)
    11 = @(0.46304099999999998, 1.467179, 1.4309719999999999)
    12 = @(0.45367499999999999, 1.1419649999999999, 1.915395)
    13 = @(0.35351399999999999, 1.394776, 1.284923)
    14 = @(0.79840100000000003, 0.84431, 1.1895690000000001)
    15 = @(0.53415400000000002, 0.81103999999999998, 0.57713199999999998)
    16 = @(0.54217000000000004, 0.54217000000000004, 0.54217000000000004)
    17 = @(0.420101, 0.420101, 0.42719400000000002)
    18 = @(0.39005699999999999, 0.41220800000000002, 0.410082)
    19 = @(0.43096000000000001, 0.430338, 0.52908299999999997)
    20 = @(0.434722, 0.434722, 0.53455900000000001)
    21 = @(0.37440699999999999, 0.39099200000000001, 0.391461)
    22 = @(0.33624700000000002, 0.33172200000000002, 0.36504900000000001)
    23 = @(0.40323900000000001, 0.35052499999999998, 0.37901499999999999)
    24 = @(1.258313, 1.258313, 1.258313)
}

foreach ($row in $results.Keys) {
    $vals = $results[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}

# Leave the selection where the author ended up after entering the data.
$ws.Range("K25").Select()
